$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.553746767090015
$ws.Range("L2").Value = 0.490485989036895

$ws.Range("B3").Value = 0.647210522905015

$ws.Range("B4").Value = 0.671354010141054
$ws.Range("C4").Value = 0.731501666372061
$ws.Range("L4").Value = 0.63442019211072

$ws.Range("B5").Value = 0.684899159160603
$ws.Range("L5").Value = 0.601821932205212

$ws.Range("B6").Value = 0.618798692282585
$ws.Range("L6").Value = 0.540183622108344
